$wb = $excel.ActiveWorkbook

# Add the new "異動名單" worksheet as the last tab
$ws = $wb.Worksheets.Add([Type]::Missing, $wb.Sheets($wb.Sheets.Count))
$ws.Name = '異動名單'

$ws.Range("A1").Value = 'with DATA as ('
$ws.Range("A2").Value = '  select '
$ws.Range("A3").Value = '   "EmpNo"'
$ws.Range("A4").Value = '  ,"AreaItem"'
$ws.Range("A5").Value = '  ,"EffectiveDate"'
$ws.Range("A6").Value = '  ,"IneffectiveDate"'
$ws.Range("A7").Value = '  ,"EmpClass"'
$ws.Range("A8").Value = '  ,ROW_NUMBER() OVER (Partition By "EmpNo"    '
$ws.Range("A9").Value = '                   '
$ws.Range("B9").Value = '    ORDER BY "EffectiveDate" Desc'
$ws.Range("A10").Value = '  '
$ws.Range("B10").Value = '                ) AS "ROWNUMBER"  '
$ws.Range("A11").Value = '  from ('
$ws.Range("A12").Value = '        select * '
$ws.Range("A13").Value = '        from "PfCoOfficer"'
$ws.Range("A14").Value = '        where "EffectiveDate" >= :startdate'
$ws.Range("A15").Value = '        union '
$ws.Range("A16").Value = '        select * '
$ws.Range("A17").Value = '        from "PfCoOfficer"'
$ws.Range("A19").Value = '       )'
$ws.Range("A20").Value = ')'
$ws.Range("A21").Value = 'select  '
$ws.Range("A22").Value = ' a."EmpNo"'
$ws.Range("A23").Value = ',a."AreaItem"'
$ws.Range("A24").Value = ',a."EffectiveDate"'
$ws.Range("A25").Value = ',a."IneffectiveDate"'
$ws.Range("A26").Value = ',a."EmpClass"'
$ws.Range("A27").Value = ',nvl(b."EmpClass",'''') as "PrevEmpClass"'
$ws.Range("A28").Value = ',case when a."IneffectiveDate" <= :enddate then ''刪除－離職'''
$ws.Range("A29").Value = '      ELSE  '' '''
$ws.Range("A30").Value = ' end as "ChangeReason"'
$ws.Range("A32").Value = 'from DATA a'
$ws.Range("A33").Value = 'left join DATA b on b."EmpNo" =  a."EmpNo" '
$ws.Range("A34").Value = '               and b.ROWNUMBER = 2'
$ws.Range("A35").Value = '               and a."IneffectiveDate" >= :enddate'
$ws.Range("A36").Value = 'where a.ROWNUMBER = 1'
$ws.Range("A37").Value = '  and case when a."IneffectiveDate" <= :enddate then 1'
$ws.Range("A38").Value = '           when a."EmpClass" <> nvl(b."EmpClass",a."EmpClass") then 2'
$ws.Range("A39").Value = '           else 0'
$ws.Range("A40").Value = '      end > 0     '
$ws.Range("A41").Value = ';'
$ws.Range("A42").Value = '參數'
$ws.Range("A18").Value = '        where "EffectiveDate" between :startdate and :enddate'
$ws.Range("B42").Value = 'startdate 季初日'
$ws.Range("B43").Value = 'enddate 季底日'

# Column B width to fit the parameter descriptions
$ws.Columns.Item(2).ColumnWidth = 34.2857142857143

# Make this the active sheet/tab, then reproduce the source file's scroll
# position and cell selection
$ws.Activate()
$ws.Range("A28").Select()
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("D41").Select()
